$d = $word.ActiveDocument

# --- 1. Refresh the footer "generated at" timestamp -------------------------
foreach ($story in $d.StoryRanges) {
    [void]$story.Find.Execute("2025-06-30 12:13Z / ", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "2025-07-02 02:48Z / ", 2)
}

# --- 2. Add the regression-test character styles (b, i, sub, sup, u) -------
# wdStyleTypeCharacter = 2

$bold = $d.Styles.Add("b", 2)
$bold.BaseStyle = "DefaultParagraphFont"
$bold.Priority = 1
$bold.QuickStyle = $true
$bold.Font.Bold = $true

$italic = $d.Styles.Add("i", 2)
$italic.BaseStyle = "DefaultParagraphFont"
$italic.Priority = 1
$italic.QuickStyle = $true
$italic.Font.Italic = $true

$sub = $d.Styles.Add("sub", 2)
$sub.BaseStyle = "DefaultParagraphFont"
$sub.Priority = 1
$sub.QuickStyle = $true
$sub.Font.Subscript = $true

$sup = $d.Styles.Add("sup", 2)
$sup.BaseStyle = "DefaultParagraphFont"
$sup.Priority = 1
$sup.QuickStyle = $true
$sup.Font.Superscript = $true

# wdUnderlineSingle = 1
$u = $d.Styles.Add("u", 2)
$u.BaseStyle = "DefaultParagraphFont"
$u.Priority = 1
$u.QuickStyle = $true
$u.Font.Underline = 1
